# DDAS Upload Template edit
# Adds PI / Institute / Sub-Investigator columns, fills in sample PI data,
# formats the header row (bold, bordered, centered) and the data row
# (centered / wrapped), widens columns, sets row height and page setup.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlCenter = -4108

# ---------------------------------------------------------------------
# 1. Header row (row 1) text, A..AC
# ---------------------------------------------------------------------
$headers = @{
    "A1"="PI Name"; "B1"="PI Medical license #"; "C1"="PI Qualification";
    "D1"="Project Number"; "E1"="Sponsor Protocol #"; "F1"="Institute Name";
    "G1"="Address"; "H1"="Country"; "I1"="Sub Investigator";
    "J1"="Sub Investigator ML#"; "K1"="SI Qualification";
    "L1"="Sub Investigator"; "M1"="Sub Investigator ML#"; "N1"="SI Qualification";
    "O1"="Sub Investigator"; "P1"="Sub Investigator ML#"; "Q1"="SI Qualification";
    "R1"="Sub Investigator"; "S1"="Sub Investigator ML#"; "T1"="SI Qualification";
    "U1"="Sub Investigator"; "V1"="Sub Investigator ML#"; "W1"="SI Qualification";
    "X1"="Sub Investigator"; "Y1"="Sub Investigator ML#"; "Z1"="SI Qualification";
    "AA1"="Sub Investigator"; "AB1"="Sub Investigator ML#"; "AC1"="SI Qualification"
}
foreach ($addr in @("A1","B1","C1","D1","E1","F1","G1","H1","I1","J1","K1","L1","M1","N1","O1","P1","Q1","R1","S1","T1","U1","V1","W1","X1","Y1","Z1","AA1","AB1","AC1")) {
    $ws.Range($addr).Value = $headers[$addr]
}

# ---------------------------------------------------------------------
# 2. Data row (row 2) sample values
# ---------------------------------------------------------------------
$ws.Range("A2").Value = "Donald Hricik"
$ws.Range("B2").Value = 35.047761000000001
$ws.Range("C2").Value = "MD"
$ws.Range("D2").Value = "0078-0609"
$ws.Range("E2").Value = "IM103116"
$instituteText = "University Hospitals Cleveland Medical Center`n11100 Euclid Avenue`nCleveland , OH 44106 / USA`n"
$ws.Range("G2").Value = $instituteText
$ws.Range("F2").ClearContents() | Out-Null

# ---------------------------------------------------------------------
# 3. Formatting
# ---------------------------------------------------------------------
# Header row: bold, thin border all around, centered
$headerRng = $ws.Range("A1:AC1")
$headerRng.Font.Bold = $true
$headerRng.Borders.LineStyle = 1
$headerRng.HorizontalAlignment = $xlCenter

# Data row first block: centered horizontally + vertically
$dataRng1 = $ws.Range("A2:R2")
$dataRng1.HorizontalAlignment = $xlCenter
$dataRng1.VerticalAlignment = $xlCenter

# Institute / Project number cells: also wrap text
$ws.Range("D2").WrapText = $true
$ws.Range("G2").WrapText = $true

# Remaining sub-investigator block: centered horizontally only
$dataRng2 = $ws.Range("S2:AC2")
$dataRng2.HorizontalAlignment = $xlCenter

# ---------------------------------------------------------------------
# 4. Row height / column widths
# ---------------------------------------------------------------------
$ws.Rows.Item(2).RowHeight = 92.25

$ws.Columns.Item(1).ColumnWidth = 17.736979166666668
$ws.Columns.Item(2).ColumnWidth = 17.736979166666668
$ws.Columns.Item(3).ColumnWidth = 23.022135416666668
$ws.Columns.Item(4).ColumnWidth = 14.022135416666666
$ws.Columns.Item(5).ColumnWidth = 23.022135416666668
$ws.Columns.Item(6).ColumnWidth = 13.592447916666666
$ws.Columns.Item(7).ColumnWidth = 31.451822916666668
$ws.Columns.Item(8).ColumnWidth = 16.022135416666668
$ws.Columns.Item(9).ColumnWidth = 25.166666666666668
$ws.Columns.Item(10).ColumnWidth = 18.736979166666668
$ws.Columns.Item(11).ColumnWidth = 17.022135416666668
$ws.Columns.Item(12).ColumnWidth = 14.592447916666666
$ws.Columns.Item(13).ColumnWidth = 18.736979166666668
$ws.Columns.Item(14).ColumnWidth = 18.451822916666668
$ws.Columns.Item(15).ColumnWidth = 14.592447916666666
$ws.Columns.Item(16).ColumnWidth = 18.736979166666668
$ws.Columns.Item(17).ColumnWidth = 17.736979166666668
$ws.Columns.Item(18).ColumnWidth = 14.592447916666666
$ws.Columns.Item(19).ColumnWidth = 18.736979166666668
$ws.Columns.Item(20).ColumnWidth = 13.736979166666666
$ws.Columns.Item(21).ColumnWidth = 14.592447916666666
$ws.Columns.Item(22).ColumnWidth = 18.736979166666668
$ws.Columns.Item(23).ColumnWidth = 13.736979166666666
$ws.Columns.Item(24).ColumnWidth = 14.592447916666666
$ws.Columns.Item(25).ColumnWidth = 18.736979166666668
$ws.Columns.Item(26).ColumnWidth = 13.736979166666666
$ws.Columns.Item(27).ColumnWidth = 14.592447916666666
$ws.Columns.Item(28).ColumnWidth = 18.736979166666668
$ws.Columns.Item(29).ColumnWidth = 13.736979166666666

# ---------------------------------------------------------------------
# 5. Page setup / selection
# ---------------------------------------------------------------------
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

$ws.Range("B3").Select() | Out-Null
